# NIT-8600024002.xlsx update
# - Adds a new "2507" period row at the top of the period-mora table
#   (periods are listed newest-first, descending), shifting the whole
#   table down by one row and dropping the oldest period ("2305") onto
#   the final (bottom-bordered) row.
# - Updates the VALOR MORA total and the Cant. Periodos counter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new data row above the first period row (row 16) ---------
$ws.Rows("16").Insert()

# The freshly inserted row is blank / default-styled; clone the formatting
# of the row right below it (the old row 16, now row 17, still carries the
# correct "interior" row style) so the new row matches the rest of the table.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the constant worker columns for the new row (same worker/value
# on every period row of this account statement). Use Value2 when reading
# back (Value's getter is unreliable on this host) - writing with Value
# works fine.
$ws.Cells.Item(16, 2).Value = $ws.Cells.Item(17, 2).Value2   # Tipo Doc Trabajador
$ws.Cells.Item(16, 3).Value = $ws.Cells.Item(17, 3).Value2   # N° Doc Trabajador
$ws.Cells.Item(16, 4).Value = $ws.Cells.Item(17, 4).Value2   # Nombre Trabajador
$ws.Cells.Item(16, 6).Value = $ws.Cells.Item(17, 6).Value2   # Valor Mora
$ws.Cells.Item(16, 7).Value = $ws.Cells.Item(17, 7).Value2   # Salario Basico

# --- 2. Re-populate the Periodo Mora column (E16:E42), newest -> oldest ---
$periods = @("2507","2506","2505","2504","2503","2502","2501", `
             "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401", `
             "2312","2311","2310","2309","2308","2307","2306","2305")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# --- 3. Update the summary figures -----------------------------------------
$ws.Range("E11").Value = 1366983      # VALOR MORA total
$ws.Range("F13").Value = 27           # Cant. Periodos
